# Delete the "Nephorology" worksheet, leaving only "Brown, et al."
# (The workbook was reorganized: the old recipe sheet "Nephorology" was
#  removed after its contents were superseded by the "Brown, et al." sheet.)

$wb = $excel.ActiveWorkbook

$excel.DisplayAlerts = $false

$nephSheet = $wb.Worksheets.Item("Nephorology")
$nephSheet.Delete()

$excel.DisplayAlerts = $true

# Make sure the remaining sheet is the active one / selected tab.
$brownSheet = $wb.Worksheets.Item("Brown, et al.")
$brownSheet.Activate()
$wb.Save()
